# Auto-update draw results: append the 2025-11-05 "Pick 3" draw as a new
# row at the bottom of the Results table (mirrors the existing rows, which
# are all stored as literal/text values, not numbers or dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out where the new row goes (row 50, right after the existing
# 49 data rows).
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$newData = @("2025-11-05", "Pick 3", "251105", "6-1-0", "2025-11-05T21:39:16.754+04:00")

# Pre-format the target row as Text so values like "2025-11-05" and
# "251105" are stored verbatim (t="str") instead of being auto-converted
# to a date serial / number, matching how every other row in the sheet
# is stored.
$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $newData[0]
$ws.Cells.Item($newRow, 2).Value = $newData[1]
$ws.Cells.Item($newRow, 3).Value = $newData[2]
$ws.Cells.Item($newRow, 4).Value = $newData[3]
$ws.Cells.Item($newRow, 5).Value = $newData[4]

# Drop back to the default "Normal" style so the new cells don't carry a
# leftover text-number-format style index, consistent with the rest of
# the sheet (none of the existing cells carry an explicit style).
$rowRange.Style = "Normal"
